{"js": "// Journal post and discussion board post assignments\n// Revises three sentences to read more smoothly (adds \"which\" clauses /\n// tense fixes) without changing their meaning.\n\nconst replacements = [\n  {\n    find:\n      \"The slum clearance policy intended to eliminate the squalor and poor living conditions but resulting in the destruction of industrial jobs and a domestic refugee crisis.\",\n    replace:\n      \"The slum clearance policy, which was intended to eliminate squalor and poor living conditions but resulted in the destruction of industrial jobs and a domestic refugee crisis.\",\n  },\n  {\n    // The sentence is split across two existing runs; only the second run's\n    // text (the part starting with \"duced...\") actually changes, so leave\n    // the first run (\"...already been re\") untouched.\n    find:\n      \"duced to dangerously low levels resulting in worse fire response and increased human fatalities.\",\n    replace:\n      \"duced to dangerously low levels, which resulted in worse fire response and increased human fatalities.\",\n  },\n  {\n    find:\n      \"utilitarian ethic because they assume the objective is to produce the greatest amount of good for the greatest amount of people through a market environment.  I believe a dialogic ethic informs cooperative, face-to-face problem solving, which requires discussion and an exchange of information and ideas to solve complex problems.\",\n    replace:\n      \"utilitarian ethic because they posit that the objective is to produce the greatest amount of good for the greatest amount of people through a market environment.  I believe a dialogic ethic informs cooperative, face-to-face problem solving, which requires discussion and an exchange of information and ideas to solve complex problems.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + find.substring(0, 60));\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Journal post and dicussion board post assignments\n# Revises three sentences to read more smoothly (adds \"which\" clauses /\n# tense fixes) without changing their meaning.\n\n$d = $word.ActiveDocument\n\n$find1 = \"The slum clearance policy intended to eliminate the squalor and poor living conditions but resulting in the destruction of industrial jobs and a domestic refugee crisis.\"\n$replace1 = \"The slum clearance policy, which was intended to eliminate squalor and poor living conditions but resulted in the destruction of industrial jobs and a domestic refugee crisis.\"\n\n# This sentence is split across two existing runs; only the second run's\n# text (starting at \"duced...\") actually changes, so search only within it\n# and leave the first run (\"...already been re\") untouched.\n$find2 = \"duced to dangerously low levels resulting in worse fire response and increased human fatalities.\"\n$replace2 = \"duced to dangerously low levels, which resulted in worse fire response and increased human fatalities.\"\n\n$find3 = \"utilitarian ethic because they assume the objective is to produce the greatest amount of good for the greatest amount of people through a market environment.  I believe a dialogic ethic informs cooperative, face-to-face problem solving, which requires discussion and an exchange of information and ideas to solve complex problems.\"\n$replace3 = \"utilitarian ethic because they posit that the objective is to produce the greatest amount of good for the greatest amount of people through a market environment.  I believe a dialogic ethic informs cooperative, face-to-face problem solving, which requires discussion and an exchange of information and ideas to solve complex problems.\"\n\n$pairs = @(\n    @($find1, $replace1),\n    @($find2, $replace2),\n    @($find3, $replace3)\n)\n\n# Locate each target phrase with Find, then assign the replacement text\n# directly on the found Range (rather than passing ReplaceWith/Replace to\n# Find.Execute) so only the matched text is touched.\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $ok = $rng.Find.Execute($pair[0])\n    if (-not $ok) {\n        throw \"Could not find target text: \" + $pair[0]\n    }\n    $rng.Text = $pair[1]\n}\n"}
